$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C6 (Ganados for Daniel Sanz) from 2 to 4; formula in G6 recalculates automatically
$ws.Range("C6").Value = 4

# Move the active selection to C7 (matches the new selection saved in the sheet view)
$ws.Range("C7").Select()
